# remove s3 from docs.
#
# Infrastructure diagram update: the "S3 stream" label is renamed to
# "PRIDE stream API" and the small protocol-label rectangles (FTP, Aspera,
# Globus) plus the two connectors glued to that corner of the diagram are
# reshaped/repositioned to match the new layout (S3 dropped from the set
# of documented transfer paths).
#
# NOTE on the literal point values below: Shape.Left/Top/Width/Height (and
# Adjustments) are exposed here as single-precision floats, the same as
# real PowerPoint's COM object model, and get truncated back to EMUs
# (1/12700 pt) on save. The decimal literals are chosen so that round-trip
# (pt -> float32 -> EMU) lands exactly on the target EMU coordinates from
# the OOXML diff rather than drifting by a fraction of an EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$root = $s.Shapes.Item(1)
$grp = $root.GroupItems

# --- Rectangle 26 (id 27) : "S3 stream" -> "PRIDE stream API" ---
# target EMU off=(5459055,3038002) ext=(742890,155311)
$rectS3 = $grp.Item(16)
$rectS3.Left   = 429.8468503937008
$rectS3.Top    = 239.21276092551182
$rectS3.Width  = 58.49527559055118
$rectS3.Height = 12.229212598425196
$rectS3.TextFrame.TextRange.Text = "PRIDE stream API"

# --- Rectangle 27 (id 28) : "FTP" - reposition only ---
# target EMU off=(6311975,3040514) ext=(222816,155311)
$rectFtp = $grp.Item(17)
$rectFtp.Left   = 497.005905511811
$rectFtp.Top    = 239.41055118110236
$rectFtp.Width  = 17.544566929133858
$rectFtp.Height = 12.229212598425196

# --- Rectangle 28 (id 29) : "Aspera" - reposition only ---
# target EMU off=(6718537,3037188) ext=(306896,155311)
$rectAspera = $grp.Item(18)
$rectAspera.Left   = 529.0186614173228
$rectAspera.Top    = 239.14866141732284
$rectAspera.Width  = 24.16504002007874
$rectAspera.Height = 12.229212598425196

# --- Rectangle 29 (id 30) : "Globus" - reposition only ---
# target EMU off=(7209179,3033520) ext=(402112,155311)
$rectGlobus = $grp.Item(19)
$rectGlobus.Left   = 567.6518897637795
$rectGlobus.Top    = 238.85984251968503
$rectGlobus.Width  = 31.66236220472441
$rectGlobus.Height = 12.229212598425196

# --- Straight Arrow Connector 32 (id 42) : bent connector glued to
#     Rectangle 26 / TextBox 19, reshaped + adj1 updated ---
# target EMU off=(5470987,3240272) ext=(406473,312554) adj1=24714 (0.24714)
$connBent = $grp.Item(21)
$connBent.Left   = 430.7863922127559
$connBent.Top    = 255.13953399905512
$connBent.Width  = 32.00574803149606
$connBent.Height = 24.610551181102363
$connBent.Adjustments.Item(1) = 0.24714

# --- Straight Arrow Connector 32 (id 50) : straight connector glued to
#     TextBox 48 (pridepy), reshaped to a vertical line ---
# target EMU off=(7091379,3260651) ext=(0,896033)
$connStraight = $grp.Item(24)
$connStraight.Left   = 558.3763122575983
$connStraight.Top    = 256.7441864033465
$connStraight.Width  = 0.0
$connStraight.Height = 70.55377952755906
